$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Alemania
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 138273
$ws.Range("C8").Value = 575
$ws.Range("D8").Value = 81800
$ws.Range("E8").Value = 52368
$ws.Range("F8").Value = 4288
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = 4105

# Row 16: Paises Bajos
$ws.Range("A16").Value = "Paises Bajos"
$ws.Range("B16").Value = 30449
$ws.Range("C16").Value = 1235
$ws.Range("D16").Value = 250
$ws.Range("E16").Value = 26740
$ws.Range("F16").Value = 1279
$ws.Range("G16").Value = 144
$ws.Range("H16").Value = 3459

# Row 17: Canada
$ws.Range("A17").Value = "Canada"
$ws.Range("B17").Value = 30106
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 9729
$ws.Range("E17").Value = 19182
$ws.Range("F17").Value = 557
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1195

# Row 21: India
$ws.Range("A21").Value = "India"
$ws.Range("B21").Value = 13835
$ws.Range("C21").Value = 405
$ws.Range("D21").Value = 1777
$ws.Range("E21").Value = 11606
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 452

# Row 32: Arabia Saudita
$ws.Range("A32").Value = "Arabia Saudita"
$ws.Range("B32").Value = 7142
$ws.Range("C32").Value = 762
$ws.Range("D32").Value = 1049
$ws.Range("E32").Value = 6006
$ws.Range("F32").Value = 71
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 87

# Row 33: Dinamarca
$ws.Range("A33").Value = "Dinamarca"
$ws.Range("B33").Value = 7073
$ws.Range("C33").Value = 194
$ws.Range("D33").Value = 3389
$ws.Range("E33").Value = 3348
$ws.Range("F33").Value = 93
$ws.Range("G33").Value = 15
$ws.Range("H33").Value = 336

# Row 34: Pakistan
$ws.Range("A34").Value = "Pakistan"
$ws.Range("B34").Value = 7025
$ws.Range("C34").Value = 106
$ws.Range("D34").Value = 1765
$ws.Range("E34").Value = 5125
$ws.Range("F34").Value = 46
$ws.Range("G34").Value = 7
$ws.Range("H34").Value = 135

# Row 35: Noruega
$ws.Range("A35").Value = "Noruega"
$ws.Range("B35").Value = 6905
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 32
$ws.Range("E35").Value = 6715
$ws.Range("F35").Value = 63
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 158

# Row 36: Australia
$ws.Range("A36").Value = "Australia"
$ws.Range("B36").Value = 6523
$ws.Range("C36").Value = 55
$ws.Range("D36").Value = 3819
$ws.Range("E36").Value = 2639
$ws.Range("F36").Value = 60
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 65

# Row 37: Chequia
$ws.Range("A37").Value = "Chequia"
$ws.Range("B37").Value = 6437
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 979
$ws.Range("E37").Value = 5288
$ws.Range("F37").Value = 82
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 170

# Row 67: Kazajistan
$ws.Range("A67").Value = "Kazajistan"
$ws.Range("B67").Value = 1498
$ws.Range("C67").Value = 96
$ws.Range("D67").Value = 306
$ws.Range("E67").Value = 1175
$ws.Range("F67").Value = 22
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 17

# Row 88: Principado de Andorra
$ws.Range("A88").Value = "Principado de Andorra"
$ws.Range("B88").Value = 696
$ws.Range("C88").Value = 23
$ws.Range("D88").Value = 191
$ws.Range("E88").Value = 470
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 35

# Row 89: Costa de Marfil
$ws.Range("A89").Value = "Costa de Marfil"
$ws.Range("B89").Value = 688
$ws.Range("C89").Value = 34
$ws.Range("D89").Value = 193
$ws.Range("E89").Value = 489
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 6

# Row 90: Letonia
$ws.Range("A90").Value = "Letonia"
$ws.Range("B90").Value = 682
$ws.Range("C90").Value = 7
$ws.Range("D90").Value = 88
$ws.Range("E90").Value = 589
$ws.Range("F90").Value = 5
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 5

# Row 119: Sri Lanka
$ws.Range("A119").Value = "Sri Lanka"
$ws.Range("B119").Value = 242
$ws.Range("C119").Value = 4
$ws.Range("D119").Value = 77
$ws.Range("E119").Value = 158
$ws.Range("F119").Value = 1
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 7

# Row 128: Jamaica
$ws.Range("A128").Value = "Jamaica"
$ws.Range("B128").Value = 143
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 25
$ws.Range("E128").Value = 113
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 5

# Row 147: Liberia
$ws.Range("A147").Value = "Liberia"
$ws.Range("B147").Value = 76
$ws.Range("C147").Value = 17
$ws.Range("D147").Value = 7
$ws.Range("E147").Value = 62
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 7

# Row 148: Barbados
$ws.Range("A148").Value = "Barbados"
$ws.Range("B148").Value = 75
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 15
$ws.Range("E148").Value = 55
$ws.Range("F148").Value = 4
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 5
